$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value, derived from the day's crypto price refresh.
$updates = @{
    "D2" = "34.218.83"
    "E2" = "  +1.10%  "
    "D3" = "1.783.99"
    "E3" = "  +0.12%  "
    "E4" = "  +0.21%  "
    "D5" = "226.44"
    "E5" = "  +0.84%  "
    "E6" = "  +0.32%  "
    "E7" = "  +0.20%  "
    "D8" = "31.94"
    "E8" = "  -0.39%  "
    "E9" = "  +0.97%  "
    "D10" = "0.0693"
    "E10" = "  +2.23%  "
    "E11" = "  +1.31%  "
    "D12" = "2.041.13"
    "E12" = "  +0.21%  "
    "E13" = "  -1.82%  "
    "D14" = "1.778.33"
    "E14" = "  -0.23%  "
    "D15" = "34.190.47"
    "E15" = "  +1.02%  "
    "E16" = "  +2.08%  "
    "E17" = "  +1.32%  "
    "D18" = "67.95"
    "E18" = "  +1.99%  "
    "B19" = "ShibaInu"
    "C19" = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
    "D19" = "0.0₃0802"
    "E19" = "  +3.65%  "
    "B20" = "BitcoinCash"
    "C20" = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
    "D20" = "247.17"
    "E20" = "  +3.48%  "
    "D21" = "11.00"
    "E21" = "  +3.99%  "
    "E22" = "  +0.20%  "
    "E23" = "  +2.12%  "
    "E24" = "  -0.64%  "
    "D25" = "162.41"
    "E25" = "  +1.20%  "
    "E26" = "  +2.37%  "
    "D27" = "16.31"
    "E27" = "  +1.41%  "
    "E28" = "  +1.66%  "
    "E29" = "  +0.26%  "
    "E30" = "  +0.69%  "
    "E31" = "  +2.01%  "
    "E32" = "  +4.46%  "
    "D33" = "3.71"
    "E33" = "  +5.51%  "
    "E34" = "  -1.11%  "
    "D35" = "1.445.31"
    "E35" = "  +4.23%  "
    "D36" = "0.654"
    "E36" = "  +2.38%  "
    "D37" = "2.42"
    "E37" = "  +7.65%  "
    "E38" = "  +3.76%  "
    "D39" = "1.05"
    "E39" = "  +0.85%  "
    "B40" = "Aave"
    "C40" = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
    "D40" = "80.45"
    "E40" = "  +2.37%  "
    "B41" = "HuobiToken"
    "C41" = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
    "D41" = "2.37"
    "E41" = "  -1.04%  "
    "D42" = "0.926"
    "E42" = "  +1.64%  "
    "E43" = "  +0.46%  "
    "E44" = "  +0.91%  "
    "E45" = "  +4.03%  "
    "D46" = "0.0510"
    "E46" = "  +0.87%  "
    "E47" = "  +0.03%  "
    "D48" = "0.0₆0136"
    "E48" = "  -4.02%  "
    "D49" = "1.942.95"
    "E49" = "  +0.24%  "
    "D50" = "104.78"
    "E50" = "  -2.80%  "
    "E51" = "  +0.14%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text format so numeric-looking values (e.g. "0.0510", "34.218.83")
    # are preserved exactly as strings instead of being parsed into numbers,
    # then restore the original (default) cell style so formatting is untouched.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
